$d = $word.ActiveDocument

# Locate the text to replace: " {cite:t}article." (the leading space that
# used to trail "...as proposed by" plus the {cite:t}article citation
# placeholder and the closing period). Using Find keeps this robust
# against any offset drift instead of relying on hard-coded character
# indices.
$found_range = $d.Content.Duplicate
$found_range.Find.ClearFormatting()
$found = $found_range.Find.Execute(" {cite:t}article.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the '{cite:t}article' citation placeholder to replace."
}

# Re-anchor on a fresh Range built from the match's Start/End - InsertXML
# needs a "live" Range (not one still owned by a Find operation) to
# correctly splice its replacement back into the document.
$target = $d.Range($found_range.Start, $found_range.End)

# Replace that span with plain-formatted runs: a space, the new citation
# text "Marrero et al. (2019)", and the closing period - none of them
# carrying the VerbatimChar run style that "article" had.
$replacementPkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p>' +
  '<w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Marrero et al. (2019)</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">.</w:t></w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($replacementPkg)

Write-Output "Replaced citation placeholder with 'Marrero et al. (2019)'."
